$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1: Title paragraph -- append " V.2" (bold, 16pt) after
# "Familiarity Review Template" and leave the "_GoBack" bookmark
# marking the end of the freshly-typed text (mirrors Word's own
# behaviour of dropping _GoBack at the last edit point).
# ---------------------------------------------------------------

# Locate "Familiarity Review Template" precisely via Find so we do
# not depend on hard-coded character offsets.
$titleFind = $d.Content.Find
$titleFind.Execute("Familiarity Review Template", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$titleEnd = $titleFind.Parent
$titleEnd.Collapse(0)

# Type the new text plus a throw-away trailing character; the extra
# character gives us a real (non-collapsed-at-paragraph-end) anchor
# to hang the "_GoBack" bookmark on, which we then delete so the
# bookmark ends up sitting, empty, right after " V.2".
$titleEnd.InsertAfter(" V.2#")

$newTextRange = $titleEnd.Duplicate
$newTextRange.MoveEnd(1, -1)
$newTextRange.Font.Bold = $true
$newTextRange.Font.Size = 16

$bookmarkAnchor = $newTextRange.Duplicate
$bookmarkAnchor.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bookmarkAnchor)

$throwAway = $newTextRange.Duplicate
$throwAway.Collapse(0)
$throwAway.MoveEnd(1, 1)
$throwAway.Delete()

# ---------------------------------------------------------------
# Change 2: Date "11/16" -> "12/1"
# ---------------------------------------------------------------
$d.Content.Find.Execute("11/16", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "12/1", 2) | Out-Null

# ---------------------------------------------------------------
# Change 3: merge the split "Summary: I" / "t's the diagram..."
# runs back into a single run (also drops the stray "_GoBack"
# bookmark that used to live between them).
# ---------------------------------------------------------------
$d.Content.Find.Execute("Summary: It", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Summary: It", 2) | Out-Null
